$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("H4").Value = 3.65
$ws.Range("I4").Value = 2.42
$ws.Range("J4").Value = 1.05
$ws.Range("K4").Value = 8.75
$ws.Range("L4").Value = 1.25
$ws.Range("M4").Value = 3.7
$ws.Range("N4").Value = 1.75
$ws.Range("O4").Value = 2.02
$ws.Range("P4").Value = 1.35
$ws.Range("Q4").Value = 3.05
$ws.Range("R4").Value = 1.62
$ws.Range("S4").Value = 2.15
$ws.Range("T4").Value = 9.5
$ws.Range("U4").Value = 15
$ws.Range("X4").Value = 22
$ws.Range("Y4").Value = 29
$ws.Range("Z4").Value = 8.75
$ws.Range("AA4").Value = 7.4
$ws.Range("AB4").Value = 14
$ws.Range("AC4").Value = 60
$ws.Range("AD4").Value = 9.25
$ws.Range("AE4").Value = 13.5
$ws.Range("AH4").Value = 20
$ws.Range("AI4").Value = 28
$ws.Range("AJ4").Value = 400
# Row 5
$ws.Range("G5").Value = 2.92
$ws.Range("I5").Value = 2.45
$ws.Range("L5").Value = 1.39
$ws.Range("P5").Value = 1.47
$ws.Range("Q5").Value = 2.57
$ws.Range("R5").Value = 1.83
$ws.Range("S5").Value = 1.87
$ws.Range("T5").Value = 8
$ws.Range("U5").Value = 15.5
$ws.Range("X5").Value = 29
$ws.Range("AA5").Value = 6.3
$ws.Range("AB5").Value = 16
$ws.Range("AC5").Value = 90
$ws.Range("AD5").Value = 7.3
$ws.Range("AE5").Value = 12.5
$ws.Range("AG5").Value = 29
$ws.Range("AI5").Value = 37
# Row 7
$ws.Range("G7").Value = 2.95
$ws.Range("I7").Value = 2.3
$ws.Range("U7").Value = 13.5
$ws.Range("W7").Value = 35
$ws.Range("X7").Value = 29
$ws.Range("Z7").Value = 7.6
$ws.Range("AD7").Value = 6.5
$ws.Range("AE7").Value = 10
# Row 9
$ws.Range("I9").Value = 2.95
$ws.Range("AH9").Value = 30
# Row 10
$ws.Range("G10").Value = 2.05
$ws.Range("I10").Value = 3.6
$ws.Range("L10").Value = 1.53
$ws.Range("M10").Value = 2.38
$ws.Range("AH10").Value = 41
# Row 11
$ws.Range("H11").Value = 3.9
$ws.Range("I11").Value = 4.55
$ws.Range("M11").Value = 3.7
$ws.Range("T11").Value = 8.75
$ws.Range("Z11").Value = 13.5
$ws.Range("AC11").Value = 55
$ws.Range("AD11").Value = 15.5
$ws.Range("AE11").Value = 28
$ws.Range("AG11").Value = 75
$ws.Range("AJ11").Value = 350
# Row 12
$ws.Range("H12").Value = 3.6
# Row 13
$ws.Range("G13").Value = 2.22
$ws.Range("I13").Value = 2.9
$ws.Range("T13").Value = 9
$ws.Range("U13").Value = 12
$ws.Range("V13").Value = 8.75
$ws.Range("X13").Value = 17
$ws.Range("Y13").Value = 24
$ws.Range("AD13").Value = 10
$ws.Range("AE13").Value = 15.5
$ws.Range("AF13").Value = 10.5
$ws.Range("AG13").Value = 35
$ws.Range("AH13").Value = 24
$ws.Range("AI13").Value = 29
# Row 19
$ws.Range("K19").Value = 12
$ws.Range("L19").Value = 1.25
$ws.Range("M19").Value = 3.75
$ws.Range("N19").Value = 1.8
$ws.Range("O19").Value = 2
$ws.Range("R19").Value = 1.67
$ws.Range("S19").Value = 2.1
$ws.Range("T19").Value = 10
$ws.Range("Y19").Value = 29
$ws.Range("Z19").Value = 12
$ws.Range("AD19").Value = 9
$ws.Range("AH19").Value = 19
$ws.Range("AI19").Value = 26
# Row 20
$ws.Range("J20").Value = 1.08
$ws.Range("K20").Value = 8
# Row 23
$ws.Range("J23").Value = 1.08
$ws.Range("K23").Value = 7.5
$ws.Range("R23").Value = 2.1
$ws.Range("S23").Value = 1.67
$ws.Range("Z23").Value = 7.5
$ws.Range("AD23").Value = 11
# Row 27
$ws.Range("G27").Value = 2.45
$ws.Range("H27").Value = 2.95
$ws.Range("I27").Value = 2.92
$ws.Range("L27").Value = 1.42
$ws.Range("M27").Value = 2.45
$ws.Range("N27").Value = 2.22
$ws.Range("O27").Value = 1.52
$ws.Range("S27").Value = 1.7
$ws.Range("T27").Value = 6.5
$ws.Range("U27").Value = 10.75
$ws.Range("V27").Value = 10
$ws.Range("W27").Value = 26
$ws.Range("X27").Value = 24
$ws.Range("Y27").Value = 40
$ws.Range("Z27").Value = 7.1
$ws.Range("AA27").Value = 5.8
$ws.Range("AB27").Value = 16
$ws.Range("AD27").Value = 7.6
$ws.Range("AE27").Value = 14
$ws.Range("AF27").Value = 10.75
$ws.Range("AG27").Value = 37
$ws.Range("AH27").Value = 28
$ws.Range("AI27").Value = 40
# Row 31
$ws.Range("G31").Value = 2.05
$ws.Range("I31").Value = 3.4
# Row 32
$ws.Range("G32").Value = 2.8
$ws.Range("H32").Value = 3.25
$ws.Range("I32").Value = 2.3
$ws.Range("K32").Value = 7.8
$ws.Range("M32").Value = 3.55
$ws.Range("O32").Value = 1.95
$ws.Range("R32").Value = 1.6
$ws.Range("T32").Value = 10.25
$ws.Range("U32").Value = 15.5
$ws.Range("V32").Value = 10
$ws.Range("W32").Value = 35
$ws.Range("Y32").Value = 27
$ws.Range("Z32").Value = 7.8
$ws.Range("AA32").Value = 6.5
$ws.Range("AB32").Value = 12
$ws.Range("AD32").Value = 9.25
$ws.Range("AG32").Value = 25
$ws.Range("AH32").Value = 17.5
$ws.Range("AJ32").Value = 300
# Row 34
$ws.Range("G34").Value = 1.38
$ws.Range("H34").Value = 4.45
$ws.Range("I34").Value = 6.7
$ws.Range("P34").Value = 1.34
$ws.Range("Q34").Value = 3
$ws.Range("T34").Value = 7
$ws.Range("AB34").Value = 21
$ws.Range("AH34").Value = 80
$ws.Range("AI34").Value = 75
# Row 36
$ws.Range("H36").Value = 3.3
$ws.Range("P36").Value = 1.42
$ws.Range("Q36").Value = 2.65
$ws.Range("R36").Value = 1.78
$ws.Range("S36").Value = 1.93
$ws.Range("T36").Value = 7.4
$ws.Range("U36").Value = 10
$ws.Range("X36").Value = 17.5
$ws.Range("Y36").Value = 28
$ws.Range("AD36").Value = 9.75
$ws.Range("AH36").Value = 28
$ws.Range("AI36").Value = 35
# Row 38
$ws.Range("H38").Value = 3.3
$ws.Range("X38").Value = 19
